$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Copy()
$ws.Range("B11:B14").PasteSpecial(-4122)

$ws.Range("B11").Value = "A"
$ws.Range("C11").Value = 919.38095238095241
$ws.Range("B12").Value = "B"
$ws.Range("C12").Value = 1063.047619047619
$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = 934.61904761904759
$ws.Range("B14").Value = "D"
$ws.Range("C14").Value = 960.19047619047615
